$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-24 Thursday" "2024-10-25 Friday"

Replace-Text "84×75=6300" "83×40=3320"
Replace-Text "48×24=1152" "40×71=2840"
Replace-Text "45×33=1485" "20×77=1540"
Replace-Text "88×17=1496" "44×94=4136"
Replace-Text "89×98=8722" "74×38=2812"
Replace-Text "89×64=5696" "92×39=3588"
Replace-Text "36×30=1080" "82×55=4510"
Replace-Text "94×22=2068" "46×25=1150"
Replace-Text "77×60=4620" "78×58=4524"
Replace-Text "89×81=7209" "25×79=1975"
Replace-Text "25×74=1850" "93×12=1116"
Replace-Text "21×88=1848" "91×68=6188"
Replace-Text "48×82=3936" "14×62=868"
Replace-Text "88×57=5016" "62×98=6076"
Replace-Text "46×14=644" "71×34=2414"
Replace-Text "62×73=4526" "34×15=510"
Replace-Text "92×67=6164" "28×40=1120"
Replace-Text "19×43=817" "57×99=5643"
Replace-Text "90×53=4770" "64×46=2944"
Replace-Text "64×68=4352" "30×66=1980"
Replace-Text "29×93=2697" "98×62=6076"
Replace-Text "74×45=3330" "47×23=1081"
Replace-Text "92×45=4140" "12×96=1152"
Replace-Text "16×12=192" "61×98=5978"
Replace-Text "56×96=5376" "27×62=1674"

Write-Output "Done"
